# Remove the affiliation lines from the two "Author" paragraphs, keeping
# only the "Name, Ph.D." line. The surname is wrapped in proofErr
# spell-check markers and split into three runs, matching the target markup:
#   <w:r><w:t>First Last-</w:t></w:r>
#   <w:proofErr w:type="spellStart"/>
#   <w:r><w:t>Surname</w:t></w:r>
#   <w:proofErr w:type="spellEnd"/>
#   <w:r><w:t>, Ph.D.</w:t></w:r>

$d = $word.ActiveDocument

function Set-AuthorNameOnly($StartsWith, $Xml) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text -like "$StartsWith*") {
            $start = $para.Range.Start
            $end = $para.Range.End - 1   # exclude the paragraph mark
            $target = $d.Range($start, $end)
            $target.InsertXML($Xml)
            return
        }
    }
}

# "Diego Mendez-Carbajo, Ph.D." + affiliation lines -> name only, surname wrapped in proofErr
$xml1 = '<w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr>' +
        '<w:r><w:t>Diego Mendez-</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Carbajo</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>, Ph.D.</w:t></w:r>' +
        '</w:p>'
Set-AuthorNameOnly "Diego Mendez-Carbajo" $xml1

# "Alejandro Dellachiesa, Ph.D." + affiliation lines -> name only, surname wrapped in proofErr
$xml2 = '<w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Alejandro </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Dellachiesa</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>, Ph.D.</w:t></w:r>' +
        '</w:p>'
Set-AuthorNameOnly "Alejandro Dellachiesa" $xml2
